# ==========================================================================
# Edit script: clean up outputs (strip stray quotes, convert space-separated
# quoted tokens into Python-list-style strings, convert "ip mask" pairs to
# CIDR notation) and extend the address-objects sheet with new columns/rows.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------------
# Sheet: interfaces
# --------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("interfaces")

$ws1.Range("G2").Value  = "['ping', 'fgfm']"
$ws1.Range("G3").Value  = "['ping', 'fgfm']"
$ws1.Range("G4").Value  = "['ping', 'https', 'fgfm', 'fabric']"
$ws1.Range("K4").Value  = "10.10.10.1/24"
$ws1.Range("L5").Value  = "To-LAN"
$ws1.Range("L6").Value  = "To-LAN"
$ws1.Range("G7").Value  = "['ping', 'https', 'ssh', 'http', 'fgfm', 'fabric']"
$ws1.Range("K7").Value  = "100.64.24.237/30"
$ws1.Range("M7").Value  = "FortiExtender"
$ws1.Range("G14").Value = "['ping', 'https', 'ssh', 'fgfm', 'fabric']"
$ws1.Range("K14").Value = "192.168.1.99/24"
$ws1.Range("G15").Value = "['ping', 'fabric']"
$ws1.Range("K15").Value = "169.254.1.1/24"
$ws1.Range("P15").Value = "['a', 'b']"
$ws1.Range("G16").Value = "['ping', 'https', 'ssh', 'snmp']"
$ws1.Range("K16").Value = "100.64.24.5/32"
$ws1.Range("G17").Value = "['ping']"
$ws1.Range("K17").Value = "100.66.11.22/30"
$ws1.Range("G18").Value = "['ping']"
$ws1.Range("K18").Value = "89.147.120.215/31"
$ws1.Range("G19").Value = "['ping']"
$ws1.Range("K19").Value = "100.66.11.26/30"
$ws1.Range("G21").Value = "['ping']"
$ws1.Range("K21").Value = "100.66.11.134/32"
$ws1.Range("G22").Value = "['ping']"
$ws1.Range("K22").Value = "100.66.11.166/32"
$ws1.Range("G23").Value = "['ping']"
$ws1.Range("K23").Value = "192.168.5.1/24"
$ws1.Range("P23").Value = "['internal1', 'internal2']"
$ws1.Range("G24").Value = "['ping']"
$ws1.Range("K24").Value = "89.147.120.213/31"

# --------------------------------------------------------------------------
# Sheet: address-objects
# New columns: comment, end_ip, start_ip inserted before associated-interface
# Many new rows added (private/reserved subnets, NetSG hosts) and the
# existing "nw-apNN" rows shifted down; subnet values converted to CIDR.
# --------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("address-objects")

# Clear out all the old data (rows 2-13) before writing the new table.
$ws2.Range("A2:E13").ClearContents()

# Give the three new header cells (F1/G1/H1) the same style as the existing
# bold/centered header cells before writing their values.
$ws2.Range("C1").Copy($ws2.Range("F1"))
$ws2.Range("D1").Copy($ws2.Range("G1"))
$ws2.Range("E1").Copy($ws2.Range("H1"))

$addrData = @(
    @("uuid","name","comment","end_ip","start_ip","associated-interface","type","subnet"),
    @("7a742462-c13a-51ec-12cb-aca6cd89acf1","all","","","","","",""),
    @("154aeb56-c6bd-51ec-c4d1-65232770be4c","FIREWALL_AUTH_PORTAL_ADDRESS","","","","","",""),
    @("154af4ac-c6bd-51ec-76c3-27c06cf279bb","FABRIC_DEVICE","IPv4 addresses of Fabric Devices.","","","","",""),
    @("154d08b4-c6bd-51ec-61c0-b9ff10463e07","SSLVPN_TUNNEL_ADDR1","","10.212.134.210","10.212.134.200","ssl.root","iprange",""),
    @("17ed18fa-c7c3-51ec-50c1-c283f73dbe8c","FCTEMS_ALL_FORTICLOUD_SERVERS","","","","","dynamic",""),
    @("9dd63474-c17e-51ec-414e-d391de278269","Private RFC1918 10.0.0.0","","","","","","10.0.0.0/8"),
    @("accc70ba-c17e-51ec-a150-34dd5295486c","Private RFC1918 172.16.0.0","","","","","","172.16.0.0/12"),
    @("bacc7c8c-c17e-51ec-1b44-78ac0eaac98d","Private RFC1918 192.168.0.0","","","","","","192.168.0.0/16"),
    @("f4c6c104-c17e-51ec-79d7-995927036c68","Private RFC6598 100.64.0.0","","","","","","100.64.0.0/10"),
    @("7193e550-c17e-51ec-5133-435ce1ae000e","NetSG AU01 Mgmt DMZ Subnet","","","","","","100.68.1.0/24"),
    @("7e302878-c17e-51ec-2aaa-add873a5d0b5","NetSG AU02 Mgmt DMZ Subnet","","","","","","100.68.2.0/24"),
    @("061f6e06-c17f-51ec-9ace-68cf37c88867","NetSG NTP Server 1","","","","","","103.74.171.50/32"),
    @("10262f7a-c17f-51ec-423e-2abd67c5a3e0","NetSG NTP Server 2","","","","","","103.99.241.50/32"),
    @("8f652b66-c17e-51ec-d179-d72dbd46e0e5","NetSG SNCK Loopback Allocation","","","","","","100.64.24.0/24"),
    @("19948d0e-c17f-51ec-208d-651c24467a49","NetSG Tacacs Server 1","","","","","","100.68.1.133/32"),
    @("2353eaf6-c17f-51ec-935d-653521fb65a3","NetSG Tacacs Server 2","","","","","","100.68.2.133/32"),
    @("ea5b02fe-e87d-51ec-45de-1058c8dc3015","nw-ap01","","","","","","192.168.7.10/32"),
    @("f55a2360-e87d-51ec-a5d9-f95c9392cb67","nw-ap02","","","","","","192.168.7.11/32"),
    @("017a06d8-e87e-51ec-ab97-93410a582ec4","nw-ap03","","","","","","192.168.7.13/32"),
    @("0ba7bf88-e87e-51ec-76cc-c7cc7b08155b","nw-ap04","","","","","","192.168.7.14/32"),
    @("17a955da-e87e-51ec-e257-f807bbfdb39b","nw-ap05","","","","","","192.168.7.15/32"),
    @("265076fe-e87e-51ec-ecc0-4471a7f569c0","nw-ap06","","","","","","192.168.7.16/32"),
    @("2e109428-e87e-51ec-6083-96f05f3296af","nw-ap07","","","","","","192.168.7.17/32")
)
$ws2.Range("A1:H24").Value = $addrData

# --------------------------------------------------------------------------
# Sheet: firewall-addrgrp
# --------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("firewall-addrgrp")

$ws3.Range("A2").Value = "['Private', 'RFC1918', '10.0.0.0', 'Private', 'RFC1918', '172.16.0.0', 'Private', 'RFC1918', '192.168.0.0', 'Private', 'RFC6598', '100.64.0.0']"
$ws3.Range("C2").Value = "Private Address Spaces"
$ws3.Range("A3").Value = "['NetSG', 'AU01', 'Mgmt', 'DMZ', 'Subnet', 'NetSG', 'AU02', 'Mgmt', 'DMZ', 'Subnet']"
$ws3.Range("C3").Value = "NetSG Mgmt DMZ Subnets"
$ws3.Range("A4").Value = "['NetSG', 'NTP', 'Server', '1', 'NetSG', 'NTP', 'Server', '2']"
$ws3.Range("C4").Value = "NetSG NTP Servers"
$ws3.Range("A5").Value = "['NetSG', 'Tacacs', 'Server', '1', 'NetSG', 'Tacacs', 'Server', '2']"
$ws3.Range("C5").Value = "NetSG Tacacs Servers"
$ws3.Range("A6").Value = "['nw-ap01', 'nw-ap02', 'nw-ap03', 'nw-ap04', 'nw-ap05', 'nw-ap06', 'nw-ap07']"
$ws3.Range("C6").Value = "Baulkham Hills Meraki APs"

# --------------------------------------------------------------------------
# Sheet: service-custom
# --------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("service-custom")

$ws4.Range("C2").Value  = "DNS"
$ws4.Range("C3").Value  = "HTTP"
$ws4.Range("C4").Value  = "HTTPS"
$ws4.Range("C5").Value  = "IMAP"
$ws4.Range("D5").Value  = "Email"
$ws4.Range("C6").Value  = "IMAPS"
$ws4.Range("D6").Value  = "Email"
$ws4.Range("C7").Value  = "LDAP"
$ws4.Range("D7").Value  = "Authentication"
$ws4.Range("C8").Value  = "DCE-RPC"
$ws4.Range("C9").Value  = "POP3"
$ws4.Range("D9").Value  = "Email"
$ws4.Range("C10").Value = "POP3S"
$ws4.Range("D10").Value = "Email"
$ws4.Range("C11").Value = "SAMBA"
$ws4.Range("C12").Value = "SMTP"
$ws4.Range("D12").Value = "Email"
$ws4.Range("C13").Value = "SMTPS"
$ws4.Range("D13").Value = "Email"
$ws4.Range("C14").Value = "KERBEROS"
$ws4.Range("D14").Value = "Authentication"
$ws4.Range("C15").Value = "LDAP_UDP"
$ws4.Range("D15").Value = "Authentication"
$ws4.Range("C16").Value = "SMB"
$ws4.Range("C17").Value = "ALL"
$ws4.Range("D17").Value = "General"
$ws4.Range("C18").Value = "ALL_ICMP"
$ws4.Range("D18").Value = "General"
$ws4.Range("C19").Value = "NTP"
$ws4.Range("C20").Value = "PING"
$ws4.Range("C21").Value = "SNMP"
$ws4.Range("C22").Value = "SSH"
$ws4.Range("C23").Value = "TRACEROUTE"
$ws4.Range("C24").Value = "webproxy"
$ws4.Range("C25").Value = "FortiGateAdminPort"
$ws4.Range("C26").Value = "TACACS+"
$ws4.Range("C27").Value = "udp-443"
$ws4.Range("C28").Value = "tcp-8080"

# New row 29
$ws4.Range("A29").Value = "7351"
$ws4.Range("C29").Value = "Meraki Ports"

# --------------------------------------------------------------------------
# Sheet: service-group
# --------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("service-group")

$ws5.Range("A1").Value = "members"

$ws5.Range("A2").Value = "['DNS', 'IMAP', 'IMAPS', 'POP3', 'POP3S', 'SMTP', 'SMTPS']"
$ws5.Range("B2").Value = "Email Access"
$ws5.Range("A3").Value = "['DNS', 'HTTP', 'HTTPS']"
$ws5.Range("B3").Value = "Web Access"
$ws5.Range("A4").Value = "['DCE-RPC', 'DNS', 'KERBEROS', 'LDAP', 'LDAP_UDP', 'SAMBA', 'SMB']"
$ws5.Range("B4").Value = "Windows AD"
$ws5.Range("A5").Value = "['DCE-RPC', 'DNS', 'HTTPS']"
$ws5.Range("B5").Value = "Exchange Server"

# --------------------------------------------------------------------------
# Sheet: firewall-policy
# --------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("firewall-policy")

$ws6.Range("B2").Value  = "['ALL_ICMP', 'FortiGateAdminPort', 'HTTPS', 'SNMP', 'SSH']"
$ws6.Range("C2").Value  = "always"
$ws6.Range("E2").Value  = "['NetSG', 'SNCK', 'Loopback', 'Allocation']"
$ws6.Range("F2").Value  = "['NetSG', 'Mgmt', 'DMZ', 'Subnets']"
$ws6.Range("G2").Value  = "['loopback0', 'internal3']"
$ws6.Range("H2").Value  = "['Zone_IPWan']"
$ws6.Range("J2").Value  = "NetSG Management Inbound"

$ws6.Range("B3").Value  = "['NTP']"
$ws6.Range("C3").Value  = "always"
$ws6.Range("E3").Value  = "['NetSG', 'NTP', 'Servers']"
$ws6.Range("F3").Value  = "['NetSG', 'SNCK', 'Loopback', 'Allocation']"
$ws6.Range("G3").Value  = "['Zone_IPWan']"
$ws6.Range("H3").Value  = "['loopback0']"
$ws6.Range("J3").Value  = "NetSG NTP Servers"

$ws6.Range("B4").Value  = "['TACACS+']"
$ws6.Range("C4").Value  = "always"
$ws6.Range("E4").Value  = "['NetSG', 'Tacacs', 'Servers']"
$ws6.Range("F4").Value  = "['NetSG', 'SNCK', 'Loopback', 'Allocation']"
$ws6.Range("G4").Value  = "['Zone_IPWan']"
$ws6.Range("H4").Value  = "['loopback0']"
$ws6.Range("J4").Value  = "NetSG TACACS Servers"

$ws6.Range("B5").Value  = "['PING', 'TRACEROUTE']"
$ws6.Range("C5").Value  = "always"
$ws6.Range("E5").Value  = "['Private', 'Address', 'Spaces']"
$ws6.Range("F5").Value  = "['Private', 'Address', 'Spaces']"
$ws6.Range("G5").Value  = "['any']"
$ws6.Range("H5").Value  = "['any']"
$ws6.Range("J5").Value  = "Trusted Internal ICMP"

$ws6.Range("B6").Value  = "['ALL']"
$ws6.Range("C6").Value  = "always"
$ws6.Range("E6").Value  = "['all']"
$ws6.Range("F6").Value  = "['all']"
$ws6.Range("G6").Value  = "['Zone_IPWan']"
$ws6.Range("H6").Value  = "['Zone_Data']"
$ws6.Range("J6").Value  = "LAN to IPWAN"

$ws6.Range("B7").Value  = "['ALL']"
$ws6.Range("C7").Value  = "always"
$ws6.Range("E7").Value  = "['all']"
$ws6.Range("F7").Value  = "['all']"
$ws6.Range("G7").Value  = "['Zone_Data']"
$ws6.Range("H7").Value  = "['Zone_IPWan']"
$ws6.Range("J7").Value  = "IPWAN to LAN"

$ws6.Range("B8").Value  = "['ALL']"
$ws6.Range("C8").Value  = "always"
$ws6.Range("E8").Value  = "['Private', 'Address', 'Spaces']"
$ws6.Range("F8").Value  = "['all']"
$ws6.Range("G8").Value  = "['virtual-wan-link']"
$ws6.Range("H8").Value  = "['any']"
$ws6.Range("J8").Value  = "Protect Internal"

$ws6.Range("B9").Value  = "['Meraki', 'Ports', 'PING', 'NTP']"
$ws6.Range("C9").Value  = "always"
$ws6.Range("E9").Value  = "['all']"
$ws6.Range("F9").Value  = "['Baulkham', 'Hills', 'Meraki', 'APs']"
$ws6.Range("G9").Value  = "['virtual-wan-link']"
$ws6.Range("H9").Value  = "['Zone_Data']"
$ws6.Range("J9").Value  = "Meraki Access Points"
$ws6.Range("M9").Value  = "This is only additional ports required for Speedtest.net."
$ws6.Range("N9").Value  = "certificate-inspection"

$ws6.Range("B10").Value = "['udp-443', 'tcp-8080']"
$ws6.Range("C10").Value = "always"
$ws6.Range("E10").Value = "['all']"
$ws6.Range("F10").Value = "['all']"
$ws6.Range("G10").Value = "['virtual-wan-link']"
$ws6.Range("H10").Value = "['Zone_Data']"
$ws6.Range("J10").Value = "Speedtest Access"
$ws6.Range("M10").Value = "This is only additional ports required for Speedtest.net."
$ws6.Range("N10").Value = "certificate-inspection"

$ws6.Range("B11").Value = "['ALL']"
$ws6.Range("C11").Value = "always"
$ws6.Range("E11").Value = "['all']"
$ws6.Range("F11").Value = "['all']"
$ws6.Range("G11").Value = "['virtual-wan-link']"
$ws6.Range("H11").Value = "['Zone_Data']"
$ws6.Range("J11").Value = "General Internet Access"
$ws6.Range("N11").Value = "certificate-inspection"

$ws6.Range("B12").Value = "['ALL']"
$ws6.Range("C12").Value = "always"
$ws6.Range("E12").Value = "['all']"
$ws6.Range("F12").Value = "['all']"
$ws6.Range("G12").Value = "['any']"
$ws6.Range("H12").Value = "['any']"
$ws6.Range("J12").Value = "Deny and Log"

# --------------------------------------------------------------------------
# Sheet: router-static
# --------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("router-static")

$ws7.Range("C2").Value  = "10.0.0.0/8"
$ws7.Range("C3").Value  = "100.64.0.0/10"
$ws7.Range("C4").Value  = "172.16.0.0/12"
$ws7.Range("C5").Value  = "192.168.0.0/16"
$ws7.Range("C7").Value  = "192.168.6.0/24"
$ws7.Range("C8").Value  = "192.168.7.0/24"
$ws7.Range("C9").Value  = "103.99.241.36/32"
$ws7.Range("C10").Value = "103.99.241.35/32"
$ws7.Range("C11").Value = "103.99.243.129/32"
$ws7.Range("C12").Value = "103.99.241.10/32"

Write-Host "Edit complete"
